$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.037447571754456
$ws.Range("B1").Value = 0.9730305075645447
$ws.Range("C1").Value = 0.9708099365234375
$ws.Range("D1").Value = 1.199261426925659
$ws.Range("E1").Value = 1.101291537284851
